$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.665.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.14%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.677.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.84%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''314.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.82%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.24%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.3893'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -3.51%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.3942'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -3.41%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''52.10'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -3.62%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''1.003'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +0.15%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''1.393'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -5.89%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.08644'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.09%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''25.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -4.47%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''7.312'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.77%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''7.772'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -4.48%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.00001314'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.48%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''1.757.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.25%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''93.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -3.45%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.07060'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.37%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''20.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.39%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''7.062'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.97%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''1.004'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.25%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''13.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -2.89%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''24.670.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.14%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.357'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.56%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''23.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.89%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.727'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -6.32%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''162.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.69%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''5.752'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -8.19%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''146.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -0.24%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''7.896'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -6.00%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''2.519'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +13.36%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.865.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -2.40%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''0.08386'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -5.59%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.03038'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -5.74%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.2822'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -1.41%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''6.870'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -5.72%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.9820'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -4.64%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.09489'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.48%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.556'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +5.71%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''10.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -3.20%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.7910'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -6.76%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''13.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.95%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''16.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -6.19%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.7131'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -4.24%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.563'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -6.31%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''4.189'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -1.33%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.08653'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +3.08%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +0.27%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.325'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -5.28%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''137.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -3.41%  '
$ws.Range("E51").Style = "Normal"
